$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # F3 (想去人数 for 丽水·龙泉ACG动漫游戏博览会): 1292 -> 1293
    $ws.Cells.Item(3, 6).Value = 1293

    # A new event (丽水·第四届HP国风动漫游戏嘉年华) needs to be inserted as the
    # new row 5, pushing the existing rows 5-7 down to rows 6-8. Shift the
    # rows down bottom-up by copying whole rows (values + formatting) so no
    # new/blended cell styles get fabricated the way Rows.Insert() would.
    $ws.Rows.Item(7).Copy()
    $ws.Rows.Item(8).PasteSpecial(-4104)
    $ws.Rows.Item(6).Copy()
    $ws.Rows.Item(7).PasteSpecial(-4104)
    $ws.Rows.Item(5).Copy()
    $ws.Rows.Item(6).PasteSpecial(-4104)
    $excel.CutCopyMode = 0

    # Row 5 (new): 丽水·第四届HP国风动漫游戏嘉年华
    $ws.Cells.Item(5, 1).Value = 4
    # Force the date-like string to stay text (otherwise Excel auto-converts
    # it to a date serial number), matching the other rows in this column.
    $ws.Cells.Item(5, 2).NumberFormat = "@"
    $ws.Cells.Item(5, 2).Value = "2024-07-27"
    $ws.Cells.Item(5, 3).Value = "丽水·第四届HP国风动漫游戏嘉年华"
    $ws.Cells.Item(5, 4).Value = "城北街798号 莱茵体育生活馆"
    $ws.Cells.Item(5, 5).Value = "2024.07.27 08:30-07.27 17:00"
    $ws.Cells.Item(5, 6).Value = 0
    $ws.Cells.Item(5, 7).Value = 50
    $ws.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87305"
    $ws.Cells.Item(5, 9).Value = "//i2.hdslb.com/bfs/openplatform/202406/YUnPOKGV1718268952725.jpeg"
    # B5 picked up the "@" text number format as its own cell style; reset it
    # back to the plain/default style (the stored value stays text either way).
    $ws.Cells.Item(1, 3).Copy()
    $ws.Cells.Item(5, 2).PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    # Fix the serial index numbers in column A for the rows pushed down.
    $ws.Cells.Item(6, 1).Value = 5
    $ws.Cells.Item(7, 1).Value = 6
    $ws.Cells.Item(8, 1).Value = 7

    # Row 8 (was row 7 before the shift): 想去人数 165 -> 166
    $ws.Cells.Item(8, 6).Value = 166

    # Row 8 is beyond the workbook's original used range, so the earlier
    # whole-row PasteSpecial didn't stamp its cell style onto it and the
    # later Value writes left it blank-styled; reapply column A's style.
    $ws.Cells.Item(4, 1).Copy()
    $ws.Cells.Item(8, 1).PasteSpecial(-4122)
    $excel.CutCopyMode = 0
}
